$wb = $excel.ActiveWorkbook

# --- Sheet "Đơn phụ phẫu 1" (2nd sheet): insert a new order row before the
#     existing "Tổng" row, then refresh the totals on the (now shifted) row.
$ws2 = $wb.Worksheets.Item(2)

# Push the current "Tổng" row (row 11) down to row 12 by inserting a blank
# row above it.
$ws2.Rows.Item(11).Insert()

# Fill in the new order data on row 11. The "Ngày thực hiện" column stores a
# plain text date string (matching the other rows), so force text formatting
# before the write to stop Excel from auto-parsing it into a date serial,
# then clear the format again so no style index lingers on the cell.
$ws2.Range("A11").Value = "HD-LUXURY"
$ws2.Range("B11").Value = 688
$ws2.Range("C11").NumberFormat = "@"
$ws2.Range("C11").Value = "08-27-2024"
$ws2.Range("C11").ClearFormats()
$ws2.Range("D11").Value = "LONG XUYÊN"
$ws2.Range("E11").Value = "C.hạnh"
$ws2.Range("F11").Value = "Cá nhân"
$ws2.Range("G11").Value = "Cắt mí"
$ws2.Range("H11").Value = "Đào Vương Anh"
$ws2.Range("I11").Value = 50000

# Update the totals row, now on row 12.
$ws2.Range("B12").Value = 10
$ws2.Range("I12").Value = 700000

# --- Sheet "Lương" (3rd sheet): refresh the payroll figures that depend on
#     the new "Đơn phụ phẫu 1" total.
$ws3 = $wb.Worksheets.Item(3)

$ws3.Range("B12").Value = 23
$ws3.Range("B13").Value = 3285714.285714286
$ws3.Range("B18").Value = 700000
$ws3.Range("B32").Value = 985714.2857142859
$ws3.Range("B34").Value = 985714.2857142859
